$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 12133.333
$ws.Range("N21").Value = -13069.333
$ws.Range("L21").Value = 12133.333
$ws.Range("J21").Value = 12133.333
$ws.Range("N23").Value = -12601.333
$ws.Range("L23").Value = 12133.333
$ws.Range("H23").Value = 12133.333
$ws.Range("J23").Value = 12133.333
$ws.Range("N34").ClearContents()
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 3450
$ws.Range("I34").Value = 3450
$ws.Range("M34").Value = -3247
$ws.Range("L34").Value = 0
$ws.Range("H34").Value = 3450
$ws.Range("N36").ClearContents()
$ws.Range("K36").Value = 3450
$ws.Range("I36").Value = 3450
$ws.Range("J36").Value = 0
$ws.Range("M36").Value = -2735
$ws.Range("L36").Value = 0
$ws.Range("H36").Value = 3450
$ws.Range("J62").Value = 3500
$ws.Range("I62").Value = 2499.7273
$ws.Range("K62").Value = 2499.7273
$ws.Range("M62").Value = -1875.7273
$ws.Range("L62").Value = 3500
$ws.Range("N62").Value = -4748
$ws.Range("H62").Value = 2653.6155
$ws.Range("I65").Value = 2499.7273
$ws.Range("L65").Value = 17500
$ws.Range("K65").Value = 12498.6365
$ws.Range("H65").Value = 2653.6155
$ws.Range("N65").Value = -23740
$ws.Range("J65").Value = 3500
$ws.Range("M65").Value = -9378.636500000001
$ws.Range("K76").Value = 4899.2
$ws.Range("M76").Value = -4584.2
$ws.Range("N76").Value = -6530.8
$ws.Range("L76").Value = 5900.8
$ws.Range("J76").Value = 5900.8
$ws.Range("I76").Value = 4899.2
$ws.Range("H76").Value = 5400
$ws.Range("M79").Value = -3807.2
$ws.Range("H79").Value = 5400
$ws.Range("L79").Value = 5900.8
$ws.Range("I79").Value = 4899.2
$ws.Range("K79").Value = 4899.2
$ws.Range("N79").Value = -8084.8
$ws.Range("J79").Value = 5900.8
$ws.Range("M116").Value = 292
$ws.Range("H116").Value = 2460
$ws.Range("I116").Value = 3150
$ws.Range("K116").Value = 3150
$ws.Range("J137").Value = 2520.3
$ws.Range("M137").Value = -1180.6155
$ws.Range("H137").Value = 1598.1945
$ws.Range("K137").Value = 3730.6155
$ws.Range("I137").Value = 1243.5385
$ws.Range("L137").Value = 7560.900000000001
$ws.Range("N137").Value = -12660.9
$ws.Range("I138").Value = 2162.65
$ws.Range("L138").Value = 11463.7023
$ws.Range("K138").Value = 6487.950000000001
$ws.Range("M138").Value = -1347.950000000001
$ws.Range("N138").Value = -21743.7023
$ws.Range("J138").Value = 3821.2341
$ws.Range("H138").Value = 3326.1343

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N6").Value = -9345.833000000001
$ws.Range("J6").Value = 8999.833000000001
$ws.Range("H6").Value = 22667.223
$ws.Range("L6").Value = 8999.833000000001
$ws.Range("M12").ClearContents()
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("J12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("K32").Value = 9925.757
$ws.Range("M32").Value = -9638.757
$ws.Range("I32").Value = 9925.757
$ws.Range("N32").Value = -114208.2
$ws.Range("L32").Value = 113634.2
$ws.Range("H32").Value = 28227.246
$ws.Range("J32").Value = 113634.2
$ws.Range("H63").Value = 2926.25
$ws.Range("I63").Value = 1505
$ws.Range("N63").Value = -4772
$ws.Range("J63").Value = 3400
$ws.Range("K63").Value = 1505
$ws.Range("L63").Value = 3400
$ws.Range("M63").Value = -819
$ws.Range("I66").Value = 1505
$ws.Range("N66").Value = -23864
$ws.Range("M66").Value = -4093
$ws.Range("K66").Value = 7525
$ws.Range("L66").Value = 17000
$ws.Range("H66").Value = 2926.25
$ws.Range("J66").Value = 3400
$ws.Range("L114").Value = 28000
$ws.Range("J114").Value = 28000
$ws.Range("H114").Value = 28000
$ws.Range("N114").Value = -36678
$ws.Range("K122").Value = 7561.200000000001
$ws.Range("N122").Value = -15511.4284
$ws.Range("H122").Value = 2939.0588
$ws.Range("L122").Value = 10611.4284
$ws.Range("M122").Value = -5111.200000000001
$ws.Range("I122").Value = 2520.4
$ws.Range("J122").Value = 3537.1428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J5").Value = 0
$ws.Range("I5").Value = 600
$ws.Range("K5").Value = 600
$ws.Range("L5").Value = 0
$ws.Range("H5").Value = 600
$ws.Range("N5").ClearContents()
$ws.Range("M5").Value = -487
$ws.Range("N105").Value = -338097.66
$ws.Range("I105").Value = 251495
$ws.Range("K105").Value = 251495
$ws.Range("H105").Value = 287113
$ws.Range("L105").Value = 334603.66
$ws.Range("J105").Value = 334603.66
$ws.Range("M105").Value = -249748

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M12").Value = -10000830
$ws.Range("K12").Value = 10001000
$ws.Range("L12").Value = 5250
$ws.Range("I12").Value = 10001000
$ws.Range("N12").Value = -5590
$ws.Range("J12").Value = 5250
$ws.Range("H12").Value = 3337166.5
$ws.Range("K16").Value = 1083.1666
$ws.Range("J16").Value = 143975.14
$ws.Range("M16").Value = -796.1666
$ws.Range("I16").Value = 1083.1666
$ws.Range("N16").Value = -144549.14
$ws.Range("H16").Value = 78025
$ws.Range("L16").Value = 143975.14
$ws.Range("H31").Value = 28262.436
$ws.Range("N31").Value = -51589.133
$ws.Range("L31").Value = 50999.133
$ws.Range("J31").Value = 50999.133
$ws.Range("K31").Value = 978.4
$ws.Range("M31").Value = -683.4
$ws.Range("I31").Value = 978.4
$ws.Range("N34").Value = -51403.133
$ws.Range("J34").Value = 50999.133
$ws.Range("K34").Value = 978.4
$ws.Range("I34").Value = 978.4
$ws.Range("M34").Value = -776.4
$ws.Range("L34").Value = 50999.133
$ws.Range("H34").Value = 28262.436
$ws.Range("I35").Value = 308.33334
$ws.Range("M35").Value = -14.33334000000002
$ws.Range("H35").Value = 3231.25
$ws.Range("K35").Value = 308.33334
$ws.Range("L35").Value = 12000
$ws.Range("J35").Value = 12000
$ws.Range("N35").Value = -12588
$ws.Range("N113").Value = -148315.14
$ws.Range("I113").Value = 1083.1666
$ws.Range("M113").Value = 1086.8334
$ws.Range("L113").Value = 143975.14
$ws.Range("H113").Value = 78025
$ws.Range("J113").Value = 143975.14
$ws.Range("K113").Value = 1083.1666
$ws.Range("H141").Value = 103414.29
$ws.Range("J141").Value = 68475
$ws.Range("L141").Value = 68475
$ws.Range("N141").Value = -78835

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K10").Value = 333
$ws.Range("L10").Value = 1952.0001
$ws.Range("H10").Value = 272.9
$ws.Range("M10").Value = -194
$ws.Range("N10").Value = -2230.0001
$ws.Range("J10").Value = 650.6667
$ws.Range("I10").Value = 111
$ws.Range("N34").Value = -4110
$ws.Range("J34").Value = 1314
$ws.Range("L34").Value = 3942
$ws.Range("H34").Value = 1202.7273
$ws.Range("J39").Value = 7450.75
$ws.Range("H39").Value = 6060.6
$ws.Range("L39").Value = 22352.25
$ws.Range("N39").Value = -22940.25
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("N54").Value = -7118
$ws.Range("J54").Value = 2000
$ws.Range("L54").Value = 6000
$ws.Range("H54").Value = 2000
$ws.Range("K55").Value = 101550
$ws.Range("N55").Value = -10284
$ws.Range("L55").Value = 9930
$ws.Range("J55").Value = 3310
$ws.Range("I55").Value = 33850
$ws.Range("H55").Value = 14762.5
$ws.Range("M55").Value = -101373
$ws.Range("I70").Value = 200982.4
$ws.Range("M70").Value = -602632.2
$ws.Range("J70").Value = 2985.7144
$ws.Range("L70").Value = 8957.143199999999
$ws.Range("H70").Value = 85484.336
$ws.Range("N70").Value = -9587.143199999999
$ws.Range("K70").Value = 602947.2
$ws.Range("H73").Value = 85484.336
$ws.Range("N73").Value = -11141.1432
$ws.Range("K73").Value = 602947.2
$ws.Range("M73").Value = -601855.2
$ws.Range("I73").Value = 200982.4
$ws.Range("J73").Value = 2985.7144
$ws.Range("L73").Value = 8957.143199999999
$ws.Range("I87").Value = 2588
$ws.Range("K87").Value = 7764
$ws.Range("N87").ClearContents()
$ws.Range("M87").Value = -6516
$ws.Range("H87").Value = 2588
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("H90").Value = 2588
$ws.Range("L90").Value = 0
$ws.Range("K90").Value = 23292
$ws.Range("N90").ClearContents()
$ws.Range("I90").Value = 2588
$ws.Range("M90").Value = -17052
$ws.Range("L131").Value = 4621.767
$ws.Range("N131").Value = -14701.767
$ws.Range("K131").Value = 2085.15
$ws.Range("I131").Value = 695.05
$ws.Range("H131").Value = 1358.7527
$ws.Range("M131").Value = 2954.85
$ws.Range("J131").Value = 1540.589

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K13").Value = 86.333336
$ws.Range("H13").Value = 86.333336
$ws.Range("M13").Value = 52.666664
$ws.Range("I13").Value = 86.333336
$ws.Range("K122").Value = 13279.9995
$ws.Range("N122").Value = -19358.0005
$ws.Range("H122").Value = 4623
$ws.Range("L122").Value = 14458.0005
$ws.Range("M122").Value = -10829.9995
$ws.Range("I122").Value = 4426.6665
$ws.Range("J122").Value = 4819.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K39").Value = 3000
$ws.Range("I39").Value = 3000
$ws.Range("H39").Value = 3000
$ws.Range("M39").Value = -2540

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K122").Value = 6876.333
$ws.Range("N122").ClearContents()
$ws.Range("H122").Value = 2292.111
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4426.333
$ws.Range("I122").Value = 2292.111
$ws.Range("J122").Value = 0
